$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
$els = $wb.Worksheets.Item("Elements")

# Row 2 (Extension root) - Constraint(s) column AJ: shorten ele-1 constraint text
$els.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$els.Rows.Item(2).AutoFit()

# Row 3 (Extension.id) - Type(s) column K: id -> string
$els.Range("K3").Value = "string" + [char]10
$els.Rows.Item(3).AutoFit()

# Row 6 (Extension.value[x]) - Definition column M: R4B -> R4 in URL
$els.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
